$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 32: new item #12 "Agregar JoptionPane en todos los botones" -> no
$ws.Range("A32").Value = 12
$ws.Range("B32").Value = "Agregar JoptionPane en todos los botones"
$ws.Range("C32").Value = "no"

# Row 33: new item #7.1 "Se agrego al interfaz" -> si, with start/end dates
$ws.Range("A33").Value = 7.1
$ws.Range("B33").Value = "Se agrego al interfaz"
$ws.Range("C33").Value = "si"
$ws.Range("D33").Value = 43798
$ws.Range("E33").Value = 43798

# Copy the date formatting (style index 1, numFmtId 14) from an existing
# date cell so D33/E33 reuse the same style instead of creating a new one.
$ws.Range("D8").Copy()
$ws.Range("D33:E33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update selection to match the final state (active cell E33)
$ws.Range("E33").Select()
